# Update cryptocurrency Price (D) and Volume(1h) (E) columns with refreshed
# quote data. Values are written with a leading apostrophe so Excel keeps
# them as literal text (matching the existing inline-string cell contents)
# instead of re-interpreting the numeric-looking strings as Number/Percent.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'297.63"
$ws.Range("E2").Formula = "'1.49%"
$ws.Range("D3").Formula = "'41.78"
$ws.Range("E3").Formula = "'3.03%"
$ws.Range("D4").Formula = "'5.004"
$ws.Range("E4").Formula = "'-0.35%"
$ws.Range("D5").Formula = "'0.07533"
$ws.Range("E5").Formula = "'2.84%"
$ws.Range("D6").Formula = "'1.607"
$ws.Range("E6").Formula = "'4.45%"
$ws.Range("D7").Formula = "'0.9186"
$ws.Range("E7").Formula = "'-1.00%"
$ws.Range("D9").Formula = "'0.1184"
$ws.Range("E9").Formula = "'1.40%"
$ws.Range("D10").Formula = "'0.1824"
$ws.Range("E10").Formula = "'4.81%"
$ws.Range("D11").Formula = "'0.09003"
$ws.Range("E11").Formula = "'3.08%"
$ws.Range("D12").Formula = "'0.04022"
$ws.Range("E12").Formula = "'-7.44%"
$ws.Range("E13").Formula = "'-0.48%"
$ws.Range("D14").Formula = "'0.001285"
$ws.Range("E14").Formula = "'1.48%"
$ws.Range("D15").Formula = "'0.005871"
$ws.Range("E15").Formula = "'-1.13%"
$ws.Range("D17").Formula = "'4.372"
$ws.Range("E17").Formula = "'2.05%"
$ws.Range("E18").Formula = "'1.23%"
$ws.Range("D19").Formula = "'8.265"
$ws.Range("E19").Formula = "'3.67%"
$ws.Range("D20").Formula = "'0.1370"
$ws.Range("E20").Formula = "'-2.13%"
$ws.Range("E21").Formula = "'17.41%"
$ws.Range("D22").Formula = "'0.04081"
$ws.Range("E22").Formula = "'3.54%"
$ws.Range("D23").Formula = "'0.001265"
$ws.Range("E23").Formula = "'0.28%"
$ws.Range("D24").Formula = "'0.003930"
$ws.Range("E24").Formula = "'3.75%"
$ws.Range("D25").Formula = "'0.0001302"
$ws.Range("E25").Formula = "'1.65%"
$ws.Range("D38").Formula = "'0.02413"
$ws.Range("E38").Formula = "'4.54%"
$ws.Range("D39").Formula = "'0.05198"
$ws.Range("E39").Formula = "'2.44%"
$ws.Range("D40").Formula = "'0.006300"
$ws.Range("E40").Formula = "'3.03%"
$ws.Range("D41").Formula = "'0.007812"
$ws.Range("E41").Formula = "'-0.55%"
$ws.Range("D42").Formula = "'0.1327"
$ws.Range("E42").Formula = "'3.05%"
$ws.Range("D43").Formula = "'0.007406"
$ws.Range("E43").Formula = "'0.67%"
$ws.Range("D44").Formula = "'0.007088"
$ws.Range("E44").Formula = "'-13.69%"
$ws.Range("D45").Formula = "'0.3254"
$ws.Range("E45").Formula = "'2.00%"
$ws.Range("D46").Formula = "'0.00006574"
$ws.Range("E46").Formula = "'4.47%"
$ws.Range("E47").Formula = "'-0.09%"
$ws.Range("D48").Formula = "'0.04563"
$ws.Range("E48").Formula = "'35.07%"
$ws.Range("D49").Formula = "'0.004200"
$ws.Range("E49").Formula = "'-0.03%"
$ws.Range("D50").Formula = "'0.00002100"
$ws.Range("E50").Formula = "'-0.09%"
$ws.Range("D51").Formula = "'0.0002000"
$ws.Range("E51").Formula = "'-0.09%"
